$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# These rows were corrupted: instead of the usual A..AB columns, only A + B were
# populated, with B holding a single stringified array literal, e.g.
#   [nil, nil, ..., "PANELS-(PANELS)", "Right side panel -(Right side panel )", "INTERIOR PANEL INK DIRTY", "", "", ""]
# Index 0 of that array lines up with column B, so index 21 -> W (Damage Area Name),
# index 22 -> X (Damage Part Name), index 23 -> Y (Damage Description), and the trailing
# three "" entries -> Z/AA/AB (Damage Component / Damage Type / Repair Type). Everything
# else in the array was nil, i.e. the remaining B..V columns should just be blank.

# B..V (except M/U, which carry the date format used by the Gate-In-Date/Di-Date columns)
# revert to ordinary blank, default-style cells.
$blankCols = @(2,3,4,5,6,7,8,9,10,11,12,14,15,16,17,18,19,20,22)
$dateCols  = @(13,21)

$rows = @(
    @{ Row = 9; W = 'PANELS-(PANELS)'; X = 'Right side panel -(Right side panel )'; Y = 'INTERIOR PANEL INK DIRTY' },
    @{ Row = 16; W = 'PANELS-(PANELS)'; X = 'Right side panel -(Right side panel )'; Y = 'ROOF PANEL CUT 12''X5''' },
    @{ Row = 25; W = 'PANELS-(PANELS)'; X = 'Right side panel -(Right side panel )'; Y = 'INTERIRO PANEL DIRTY' },
    @{ Row = 27; W = 'DOORS-(D)'; X = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'; Y = 'LEFT DOOR LOCKING BAR BANT 01 PES' },
    @{ Row = 29; W = 'PANELS-(PANELS)'; X = 'Right side panel -(Right side panel )'; Y = 'INTERIRO PANEL DIRTY' },
    @{ Row = 32; W = 'DOORS-(D)'; X = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'; Y = 'LEFT DOOR LOCKING BAR BANT 01 PES' },
    @{ Row = 35; W = 'DOORS-(D)'; X = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'; Y = 'LEFT DOOR LOCKING BAR BANT 01 PES' },
    @{ Row = 38; W = 'DOORS-(D)'; X = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'; Y = 'BOTH DOOR LOCKING BAR BENT 03 PCS' },
    @{ Row = 41; W = 'FLOORS-(F)'; X = 'Threshold plate-(Threshold plate)'; Y = 'FLOOR BOARD NAILS  FITTING VERYES PLS' },
    @{ Row = 42; W = 'DOORS-(D)'; X = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'; Y = 'BOTH DOOR PANEL RUSTED &  CORROSION' },
    @{ Row = 44; W = 'FLOORS-(F)'; X = 'Threshold plate-(Threshold plate)'; Y = 'FLOOR BOARD NAILS FITTING 04 PES' }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Clear the bogus stringified-array text out of column B.
    $ws.Cells.Item($r, 2).Value = ""

    # Re-touch (without changing) the format of the other "nil" columns so they persist
    # as real, explicitly-blank cells instead of staying absent from the sheet.
    foreach ($c in $blankCols) {
        $ws.Cells.Item($r, $c).Font.Name = $ws.Cells.Item($r, $c).Font.Name
    }
    foreach ($c in $dateCols) {
        $ws.Cells.Item($r, $c).NumberFormat = $ws.Cells.Item(6, $c).NumberFormat
    }

    # W/X/Y: Damage Area Name / Damage Part Name / Damage Description.
    $ws.Cells.Item($r, 23).Value = $item.W
    $ws.Cells.Item($r, 24).Value = $item.X
    $ws.Cells.Item($r, 25).Value = $item.Y

    # Z/AA/AB: Damage Component / Damage Type / Repair Type -- empty strings, not blanks.
    $ws.Cells.Item($r, 26).Formula = '=""'
    $ws.Cells.Item($r, 27).Formula = '=""'
    $ws.Cells.Item($r, 28).Formula = '=""'
}

# Column B was sized to fit the long stringified array literal; with that text gone it
# shrinks back to fit the short numeric Id values again. Column X (Damage Part Name) grows
# because several of the newly-restored rows carry some of the longest strings seen in that
# column ("Door stiffeners hinges side edge.-(...)").
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143
